$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old sub-header row (Hiver/Eté/Année labels); the power-plant
# data rows below shift up by one.
$ws.Rows.Item(2).Delete()

# Create a temporary named style that matches the font used for the data
# cells (9pt Arial) but does not force a number format, so the resulting
# cellXf only carries applyFont (matching the new header style).
$headerStyle = $wb.Styles.Add("TempHeaderStyle")
$headerStyle.Font.Size = 9
$headerStyle.IncludeNumber = $false

# New leading identifier / metadata columns.
$ws.Cells.Item(1,1).ClearFormats()
$ws.Cells.Item(1,2).ClearFormats()
$ws.Cells.Item(1,3).ClearFormats()
$ws.Cells.Item(1,4).ClearFormats()
$ws.Cells.Item(1,5).ClearFormats()

$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"

# Re-labelled measurement columns, styled with the new header style.
$ws.Range("F1:K1").Style = "TempHeaderStyle"

$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

# Drop the temporary named style - this leaves the generated cellXf behind
# (referencing the builtin Normal cellStyleXf, xfId 0) without re-adding a
# bespoke named cell style entry, and without forcing applyNumberFormat.
$headerStyle.Delete()

# Re-assert the measurement values for the shifted data rows so the stored
# doubles keep full precision (row-shift recompute can otherwise leave
# floating point noise in the serialized values).
$ws.Cells.Item(2,6).Value = 0.6
$ws.Cells.Item(2,7).Value = 0.09
$ws.Cells.Item(2,8).Value = 0.09
$ws.Cells.Item(2,9).Value = 0.2
$ws.Cells.Item(2,10).Value = 0.23
$ws.Cells.Item(2,11).Value = 0.42

$ws.Cells.Item(3,6).Value = 18
$ws.Cells.Item(3,7).Value = 8.11
$ws.Cells.Item(3,8).Value = 7.51
$ws.Cells.Item(3,9).Value = 7.66
$ws.Cells.Item(3,10).Value = 9.11
$ws.Cells.Item(3,11).Value = 16.77

$ws.Cells.Item(4,6).Value = 0.45
$ws.Cells.Item(4,7).Value = 0.72
$ws.Cells.Item(4,8).Value = 0.6
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 2

$ws.Cells.Item(5,6).Value = 0.62
$ws.Cells.Item(5,7).Value = 0.38
$ws.Cells.Item(5,8).Value = 0.33
$ws.Cells.Item(5,9).Value = 0.8
$ws.Cells.Item(5,10).Value = 1.2
$ws.Cells.Item(5,11).Value = 2

$ws.Cells.Item(6,6).Value = 2.3
$ws.Cells.Item(6,7).Value = 0.36
$ws.Cells.Item(6,8).Value = 0.36
$ws.Cells.Item(6,9).Value = 0.88
$ws.Cells.Item(6,10).Value = 0.91
$ws.Cells.Item(6,11).Value = 1.79

$ws.Range("A2:K2").Select()
